$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Marvin (sheet1): add row 31
# -------------------------------------------------------------------------
$wsMarvin = $wb.Worksheets.Item("Marvin")
$wsMarvin.Range("B30").Copy()
$wsMarvin.Range("B31").PasteSpecial(-4122)  # xlPasteFormats
$wsMarvin.Range("A31").Value = "KBS b les"
$wsMarvin.Range("B31").Value = 45280
$wsMarvin.Range("C31").Value = 120
$wsMarvin.Range("D31").Value = "les"

# -------------------------------------------------------------------------
# Demi (sheet2): add row 33
# -------------------------------------------------------------------------
$wsDemi = $wb.Worksheets.Item("Demi")
$wsDemi.Range("B32").Copy()
$wsDemi.Range("B33").PasteSpecial(-4122)
$wsDemi.Range("A33").Value = "KBS b les"
$wsDemi.Range("B33").Value = 45280
$wsDemi.Range("C33").Value = 120
$wsDemi.Range("D33").Value = "les"

# -------------------------------------------------------------------------
# Lucas (sheet3): add row 28
# -------------------------------------------------------------------------
$wsLucas = $wb.Worksheets.Item("Lucas")
$wsLucas.Range("B27").Copy()
$wsLucas.Range("B28").PasteSpecial(-4122)
$wsLucas.Range("A28").Value = "KBS b les"
$wsLucas.Range("B28").Value = 45280
$wsLucas.Range("C28").Value = 120
$wsLucas.Range("D28").Value = "les"

# -------------------------------------------------------------------------
# Luuk (sheet4): add row 34
# -------------------------------------------------------------------------
$wsLuuk = $wb.Worksheets.Item("Luuk")
$wsLuuk.Range("B33").Copy()
$wsLuuk.Range("B34").PasteSpecial(-4122)
$wsLuuk.Range("A34").Value = "KBS b les"
$wsLuuk.Range("B34").Value = 45280
$wsLuuk.Range("C34").Value = 120
$wsLuuk.Range("D34").Value = "les"

# -------------------------------------------------------------------------
# Jochem (sheet5): add row 40 (B40 only, formatted empty date cell)
# -------------------------------------------------------------------------
$wsJochem = $wb.Worksheets.Item("Jochem")
$wsJochem.Range("B38").Copy()
$wsJochem.Range("B40").PasteSpecial(-4122)

# -------------------------------------------------------------------------
# Selections / active sheet to match final saved state
# -------------------------------------------------------------------------
$wsDemi.Activate()
$wsDemi.Range("A33").Select()

$wsLucas.Activate()
$wsLucas.Range("A28").Select()

$wsLuuk.Activate()
$wsLuuk.Range("A34").Select()

$wsJochem.Activate()
$wsJochem.Range("D40").Select()

$wsMarvin.Activate()
$wsMarvin.Range("A31").Select()
